$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "540.76").
# Force it to stay text (matching the source data, which is all inline
# strings) by switching to a text number format before assigning the
# value, then restore the default "Normal" style so no stray formatting
# is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.241.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.535.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.562.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.981.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.123.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.555.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.991"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0798"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "162.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "304.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.841"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0524"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.07%  "
